$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two population figures that changed for this release
$ws.Range("C4").Value = 2891000
$ws.Range("C12").Value = 2800000

# Sort the table (header row included) by Population, descending
$dataRange = $ws.Range("A1:D13")
$dataRange.Sort($ws.Range("C1"), 2)

# Turn on AutoFilter for the table
$dataRange.AutoFilter()

# Excel records the filter range as a hidden, sheet-scoped defined name
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$13")
$filterName.Visible = $false
